# Auto-assisted edit script: apply skills.xlsx diff via Excel COM object model
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: 'specifically_assigned' -> 'game_class_id' (col M) ---
$ws.Range("M1").Value = 'game_class_id'

# --- Rows 2-4: max_level 100 -> 999, skill_bonus_per_level 0.01 -> 0.001 ---
$ws.Range("D2").Value = 999
$ws.Range("L2").Value = 0.001
$ws.Range("D3").Value = 999
$ws.Range("L3").Value = 0.001
$ws.Range("D4").Value = 999
$ws.Range("L4").Value = 0.001

# --- Insert 6 new rows before row 11 (keeps Disenchanting/Alchemy, now at 17-18) ---
$ws.Rows.Item(11).Resize(6).Insert()

# --- Row 11: Astral Magics (was: Disenchanting data, now shifted to 17; row 11 freshly populated) ---
# Row 11: Astral Magics
$ws.Range("A11").Value = 'Astral Magics'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = @'
Increases spell damage over time. The higher the level the more spell damage your magics will do. This skill is only available to heritics.
This bonus is only applied to your spell's that do damage. If you have none equipped, no bonus will be applied.
'@
$ws.Range("D11").Value = 999
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.001
$ws.Range("M11").Value = 'Heretic'
$ws.Range("N11").Value = 0

# Row 12: Celestial Prayer
$ws.Range("A12").Value = 'Celestial Prayer'
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = @'
As you level this skill over time, your healing spells will increase over time doing more and more healing.
This bonus is only applied as long as you have one healing spell equipped.
'@
$ws.Range("D12").Value = 999
$ws.Range("F12").Value = 0.001
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.001
$ws.Range("M12").Value = 'Prophet'
$ws.Range("N12").Value = 0

# Row 13: Soldiers Strength
$ws.Range("A13").Value = 'Soldiers Strength'
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = @'
This skill is only applied to fighters. Training this skill will increase your attack and defence as long as you have either a shield (for the defence bonus) or a weapon (for the attack bonus) having both will give you both bonuses, have a weapon (and no shield) or dual wielding, will only give you attack bonus.
Having double shields will only give you the defence bonus. Having double of either will not mean you get double the bonus, you will only receive the skill bonus(es) once for either type.
'@
$ws.Range("D13").Value = 999
$ws.Range("E13").Value = 0.001
$ws.Range("G13").Value = 0.001
$ws.Range("K13").Value = 1
$ws.Range("M13").Value = 'Fighter'
$ws.Range("N13").Value = 0

# Row 14: Shadow Dance
$ws.Range("A14").Value = 'Shadow Dance'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = @'
This skill is only applied to thieves and requires the thief to be dual wielding weapons to apply it's bonus. While one might assume thieves are only great with daggers and bows, in this world - a thief can duel wield any set of weapons.
While duel wielding, this skill will increase your damage over time.
This will also decrease attack time by ~5% at max level.
'@
$ws.Range("D14").Value = 999
$ws.Range("E14").Value = 0.001
$ws.Range("H14").Value = 0.005
$ws.Range("K14").Value = 1
$ws.Range("M14").Value = 'Thief'
$ws.Range("N14").Value = 0

# Row 15: Blood Lust
$ws.Range("A15").Value = 'Blood Lust'
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 'Only applies to vampires. Unlike other class specific skills vampires do not need any specific equipment for these bonuses to apply to their damage and their healing modifiers.'
$ws.Range("D15").Value = 999
$ws.Range("E15").Value = 0.001
$ws.Range("F15").Value = 0.001
$ws.Range("K15").Value = 1
$ws.Range("M15").Value = 'Vampire'
$ws.Range("N15").Value = 0

# Row 16: Natures Insight
$ws.Range("A16").Value = 'Natures Insight'
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 'This skill is only used by rangers. Increasing this skill will apply to attack and fight time out (~5% bonus at max level). Similar to thieves with the exception that these bonuses will only apply if you have a bow equipped.'
$ws.Range("D16").Value = 999
$ws.Range("E16").Value = 0.001
$ws.Range("H16").Value = 0.005
$ws.Range("K16").Value = 1
$ws.Range("M16").Value = 'Ranger'
$ws.Range("N16").Value = 0

# --- Row 17 (was row 11, Disenchanting): update description text + is_locked flag ---
$ws.Range("C17").Value = @'
Allows you to click disenchant to destroy items and get gold dust. The higher the level the more successful at getting the gold dust.
Destroying items only gets you between 1-25 Gold Dust, guaranteed, while disenchanting can get you between 1-150 or 1 (Gold Dust) if you fail to disenchant.
Gold dust is used for crafting special items with Alchemy and for conjuring celestial entities.
'@
$ws.Range("N17").Value = 1

# --- Row 18 (was row 12, Alchemy): update description text ---
$ws.Range("C18").Value = @'
Alchemy is a skill that is locked behind a quest. Once unlocked you can click Craft/Enchant to then click a new option called: Alchemy.
Alchemy requires the use of both Gold Dust (you get from disenchanting and destroying) and Shards (you get from killing (you have to be the one that kills) Celestial Entities). Once you have enough you can start crafting items that either give you limited time boons or items that do damage to kingdoms (which you can only use when attacking a kingdom).
Boons can be used from your inventory. You can use a maximum of ten boon at one time.
'@

# --- Wrap text + row heights for long-description rows ---
# Union so WrapText is applied in a single style-creating operation
$wrapRange = $excel.Union($ws.Range("C11:C14"), $ws.Range("C17:C18"))
$wrapRange.WrapText = $true

$ws.Rows.Item(11).RowHeight = 48
$ws.Rows.Item(12).RowHeight = 48
$ws.Rows.Item(13).RowHeight = 48
$ws.Rows.Item(14).RowHeight = 80
$ws.Rows.Item(17).RowHeight = 80
$ws.Rows.Item(18).RowHeight = 80

# --- Column C width update (308 -> 417) ---
$ws.Columns.Item(3).ColumnWidth = 417

# --- Selection / view tweaks to mirror the saved file's last-known state ---
$ws.Range("C18").Select()

